$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price / volume data per the latest GitHub Actions scrape run.
# Column D holds formatted price strings and column E holds formatted 1h
# volume-change percentage strings. For cells whose price text looks like a
# plain number, force text formatting first so Excel keeps the exact string
# (e.g. "244.85") instead of auto-converting it to a floating point number.
# Rows 20/21 and 33/34 also swap their B (Coin) and C (Link) values because
# two coins swapped rank positions in the source ranking.

$ws.Range("D2").Value = '36.329.29'
$ws.Range("E2").Value = '  -3.23%  '
$ws.Range("D3").Value = '1.975.38'
$ws.Range("E3").Value = '  -4.11%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.85'
$ws.Range("E5").Value = '  -3.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  -5.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.66'
$ws.Range("E7").Value = '  -12.03%  '
$ws.Range("E9").Value = '  -6.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '56.70'
$ws.Range("E10").Value = '  -6.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0838'
$ws.Range("E11").Value = '  +6.73%  '
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.00'
$ws.Range("E13").Value = '  -4.18%  '
$ws.Range("E14").Value = '  -9.20%  '
$ws.Range("D15").Value = '2.264.61'
$ws.Range("E15").Value = '  -4.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.86'
$ws.Range("E16").Value = '  -8.87%  '
$ws.Range("E17").Value = '  -6.15%  '
$ws.Range("D18").Value = '1.973.03'
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("D19").Value = '36.195.35'
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0881'
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.25'
$ws.Range("E21").Value = '  -4.88%  '
$ws.Range("E22").Value = '  -5.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.53'
$ws.Range("E23").Value = '  -3.52%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("E25").Value = '  -5.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.32'
$ws.Range("E26").Value = '  -5.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.79'
$ws.Range("E27").Value = '  -3.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.18'
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.74'
$ws.Range("E30").Value = '  -2.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.119'
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("E32").Value = '  -4.39%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0681'
$ws.Range("E33").Value = '  +6.95%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.85'
$ws.Range("E34").Value = '  -7.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.39'
$ws.Range("E35").Value = '  -6.94%  '
$ws.Range("E36").Value = '  -3.05%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("E39").Value = '  -9.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("E40").Value = '  -7.36%  '
$ws.Range("E41").Value = '  -5.29%  '
$ws.Range("E42").Value = '  -7.79%  '
$ws.Range("E43").Value = '  -5.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0213'
$ws.Range("E44").Value = '  -3.80%  '
$ws.Range("E45").Value = '  -7.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.39'
$ws.Range("E46").Value = '  -6.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.00'
$ws.Range("E47").Value = '  -12.07%  '
$ws.Range("D48").Value = '1.359.35'
$ws.Range("E48").Value = '  -4.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.42'
$ws.Range("E49").Value = '  -7.62%  '
$ws.Range("E50").Value = '  -4.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.70'
$ws.Range("E51").Value = '  -8.59%  '
